# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    # Force the cell to remain plain text so Excel does not reinterpret
    # numeric-looking strings (e.g. "553.26") as floating point numbers.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "58.337.17"
$ws.Range("E2").Value = "  -8.00%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.903.51"
$ws.Range("E3").Value = "  -8.50%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "553.26"
$ws.Range("E5").Value = "  -8.01%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "120.74"
$ws.Range("E6").Value = "  -11.22%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.29%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Range("D8") "2.899.14"
$ws.Range("E8").Value = "  -8.58%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -4.04%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.126"
$ws.Range("E10").Value = "  -11.68%  "

# Row 11 - Toncoin
Set-TextValue $ws.Range("D11") "4.82"
$ws.Range("E11").Value = "  -10.25%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -5.32%  "

# Row 13 - ShibaInu
Set-TextValue $ws.Range("D13") "0.0000212"
$ws.Range("E13").Value = "  -11.45%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "31.38"
$ws.Range("E14").Value = "  -10.16%  "

# Row 15 - TRON
Set-TextValue $ws.Range("D15") "0.119"
$ws.Range("E15").Value = "  -1.03%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D16") "3.391.15"
$ws.Range("E16").Value = "  -8.23%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.911.54"
$ws.Range("E17").Value = "  -8.24%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "58.592.14"
$ws.Range("E18").Value = "  -7.69%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "6.39"
$ws.Range("E19").Value = "  -3.12%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "415.01"
$ws.Range("E20").Value = "  -10.07%  "

# Row 21 - Chainlink
Set-TextValue $ws.Range("D21") "12.80"
$ws.Range("E21").Value = "  -8.47%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -6.71%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "6.80"
$ws.Range("E23").Value = "  -11.26%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D24") "12.55"
$ws.Range("E24").Value = "  -5.86%  "

# Row 25
Set-TextValue $ws.Range("D25") "77.10"
$ws.Range("E25").Value = "  -7.39%  "

# Row 26
Set-TextValue $ws.Range("D26") "1.00"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.27%  "

# Row 28 - PancakeSwap
Set-TextValue $ws.Range("D28") "2.45"
$ws.Range("E28").Value = "  -9.08%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  -9.31%  "

# Row 30 - RenderToken
Set-TextValue $ws.Range("D30") "6.98"
$ws.Range("E30").Value = "  -9.40%  "

# Row 31 - now EthereumClassic (was NEARProtocol)
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D31") "24.51"
$ws.Range("E31").Value = "  -9.84%  "

# Row 32 - now NEARProtocol (was EthereumClassic)
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D32") "5.92"
$ws.Range("E32").Value = "  -12.80%  "

# Row 33 - Hedera
Set-TextValue $ws.Range("D33") "0.0937"
$ws.Range("E33").Value = "  -6.75%  "

# Row 34 - Filecoin
Set-TextValue $ws.Range("D34") "5.39"
$ws.Range("E34").Value = "  -8.94%  "

# Row 35 - now Mantle (was OKB)
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D35") "0.899"
$ws.Range("E35").Value = "  -12.20%  "

# Row 36 - now OKB (was Mantle)
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D36") "48.88"
$ws.Range("E36").Value = "  -4.70%  "

# Row 37 - Stacks
Set-TextValue $ws.Range("D37") "1.96"
$ws.Range("E37").Value = "  -19.37%  "

# Row 38 - Cosmos
Set-TextValue $ws.Range("D38") "8.30"
$ws.Range("E38").Value = "  +1.95%  "

# Row 39 - PEPE
Set-TextValue $ws.Range("D39") "0.0₃0615"
$ws.Range("E39").Value = "  -16.43%  "

# Row 40 - VeChain
Set-TextValue $ws.Range("D40") "0.0343"
$ws.Range("E40").Value = "  -12.20%  "

# Row 41 - Kaspa
Set-TextValue $ws.Range("D41") "0.105"
$ws.Range("E41").Value = "  -6.86%  "

# Row 42 - Maker
Set-TextValue $ws.Range("D42") "2.618.08"
$ws.Range("E42").Value = "  -6.54%  "

# Row 43 - now dogwifhat (was Bittensor)
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D43") "2.35"
$ws.Range("E43").Value = "  -10.87%  "

# Row 44 - now Bittensor (was dogwifhat)
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D44") "350.54"
$ws.Range("E44").Value = "  -10.47%  "

# Row 46 - Monero
Set-TextValue $ws.Range("D46") "118.37"
$ws.Range("E46").Value = "  -6.00%  "

# Row 47 - TheGraph
$ws.Range("E47").Value = "  -9.79%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -5.52%  "

# Row 49 - Fetch.AI
Set-TextValue $ws.Range("D49") "1.90"
$ws.Range("E49").Value = "  -10.38%  "

# Row 50 - InjectiveProtocol
Set-TextValue $ws.Range("D50") "22.58"
$ws.Range("E50").Value = "  -10.22%  "

# Row 51 - ThetaToken
Set-TextValue $ws.Range("D51") "1.94"
$ws.Range("E51").Value = "  -10.73%  "
